$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 14 de Octubre de 2020 a las 00:56"

# Row 4
$ws.Cells.Item(4, 2).Value = 8084985
$ws.Cells.Item(4, 3).Value = 46409
$ws.Cells.Item(4, 4).Value = 5219447
$ws.Cells.Item(4, 5).Value = 2644769
$ws.Cells.Item(4, 7).Value = 751
$ws.Cells.Item(4, 8).Value = 220769

# Row 5
$ws.Cells.Item(5, 2).Value = 7237082
$ws.Cells.Item(5, 3).Value = 63517
$ws.Cells.Item(5, 4).Value = 6298695
$ws.Cells.Item(5, 5).Value = 827770

# Row 6
$ws.Cells.Item(6, 2).Value = 5113628
$ws.Cells.Item(6, 3).Value = 10220
$ws.Cells.Item(6, 4).Value = 4526975
$ws.Cells.Item(6, 5).Value = 435655
$ws.Cells.Item(6, 7).Value = 289
$ws.Cells.Item(6, 8).Value = 150998

# Row 9
$ws.Cells.Item(9, 2).Value = 924098
$ws.Cells.Item(9, 3).Value = 5015
$ws.Cells.Item(9, 4).Value = 806703
$ws.Cells.Item(9, 5).Value = 89254
$ws.Cells.Item(9, 7).Value = 156
$ws.Cells.Item(9, 8).Value = 28141

# Row 25
$ws.Cells.Item(25, 2).Value = 335679
$ws.Cells.Item(25, 3).Value = 4585
$ws.Cells.Item(25, 5).Value = 46839

# Row 30
$ws.Cells.Item(30, 2).Value = 186332
$ws.Cells.Item(30, 3).Value = 3493
$ws.Cells.Item(30, 4).Value = 157014
$ws.Cells.Item(30, 5).Value = 19669
$ws.Cells.Item(30, 7).Value = 22
$ws.Cells.Item(30, 8).Value = 9649

# Row 34
$ws.Cells.Item(34, 2).Value = 148171
$ws.Cells.Item(34, 3).Value = 856
$ws.Cells.Item(34, 5).Value = 7802
$ws.Cells.Item(34, 7).Value = 17
$ws.Cells.Item(34, 8).Value = 12235

# Row 46
$ws.Cells.Item(46, 2).Value = 104787
$ws.Cells.Item(46, 3).Value = 139
$ws.Cells.Item(46, 4).Value = 97841
$ws.Cells.Item(46, 5).Value = 875
$ws.Cells.Item(46, 7).Value = 9
$ws.Cells.Item(46, 8).Value = 6071

# Row 57
$ws.Cells.Item(57, 2).Value = 76272
$ws.Cells.Item(57, 3).Value = 324
$ws.Cells.Item(57, 4).Value = 72164
$ws.Cells.Item(57, 5).Value = 3823
$ws.Cells.Item(57, 7).Value = 5
$ws.Cells.Item(57, 8).Value = 285

# Row 58
$ws.Cells.Item(58, 5).Value = 13978
$ws.Cells.Item(58, 7).Value = 4
$ws.Cells.Item(58, 8).Value = 2103

# Row 61
$ws.Cells.Item(61, 2).Value = 60655
$ws.Cells.Item(61, 3).Value = 225
$ws.Cells.Item(61, 4).Value = 52006
$ws.Cells.Item(61, 5).Value = 7533
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 8).Value = 1116

# Row 84
$ws.Cells.Item(84, 2).Value = 27317
$ws.Cells.Item(84, 3).Value = 31
$ws.Cells.Item(84, 5).Value = 1381

# Row 85
$ws.Cells.Item(85, 2).Value = 25774
$ws.Cells.Item(85, 3).Value = 785
$ws.Cells.Item(85, 4).Value = 16139
$ws.Cells.Item(85, 5).Value = 8712
$ws.Cells.Item(85, 7).Value = 8
$ws.Cells.Item(85, 8).Value = 923

# Row 95
$ws.Cells.Item(95, 1).Value = "Noruega"
$ws.Cells.Item(95, 2).Value = 15793
$ws.Cells.Item(95, 3).Value = 154
$ws.Cells.Item(95, 4).Value = 11863
$ws.Cells.Item(95, 5).Value = 3653
$ws.Cells.Item(95, 7).Value = 1
$ws.Cells.Item(95, 8).Value = 277

# Row 96
$ws.Cells.Item(96, 1).Value = "Albania"
$ws.Cells.Item(96, 2).Value = 15752
$ws.Cells.Item(96, 3).Value = 182
$ws.Cells.Item(96, 4).Value = 9675
$ws.Cells.Item(96, 5).Value = 5648
$ws.Cells.Item(96, 7).Value = 5
$ws.Cells.Item(96, 8).Value = 429

# Row 97
$ws.Cells.Item(97, 2).Value = 15587
$ws.Cells.Item(97, 3).Value = 38
$ws.Cells.Item(97, 4).Value = 14783
$ws.Cells.Item(97, 5).Value = 459

# Row 113
$ws.Cells.Item(113, 2).Value = 8887
$ws.Cells.Item(113, 3).Value = 5
$ws.Cells.Item(113, 4).Value = 7140
$ws.Cells.Item(113, 5).Value = 1517

# Row 115
$ws.Cells.Item(115, 2).Value = 8036
$ws.Cells.Item(115, 3).Value = 15
$ws.Cells.Item(115, 4).Value = 7632
$ws.Cells.Item(115, 5).Value = 174

# Row 117
$ws.Cells.Item(117, 2).Value = 7565
$ws.Cells.Item(117, 3).Value = 11
$ws.Cells.Item(117, 4).Value = 7301
$ws.Cells.Item(117, 5).Value = 101

# Row 119
$ws.Cells.Item(119, 1).Value = "Guadalupe"
$ws.Cells.Item(119, 2).Value = 6908
$ws.Cells.Item(119, 3).Value = 425
$ws.Cells.Item(119, 4).Value = 2199
$ws.Cells.Item(119, 5).Value = 4613
$ws.Cells.Item(119, 7).Value = 19
$ws.Cells.Item(119, 8).Value = 96

# Row 120
$ws.Cells.Item(120, 1).Value = "Angola"
$ws.Cells.Item(120, 2).Value = 6680
$ws.Cells.Item(120, 3).Value = 192
$ws.Cells.Item(120, 4).Value = 2761
$ws.Cells.Item(120, 5).Value = 3697
$ws.Cells.Item(120, 7).Value = 3
$ws.Cells.Item(120, 8).Value = 222

# Row 126
$ws.Cells.Item(126, 2).Value = 5353
$ws.Cells.Item(126, 3).Value = 89
$ws.Cells.Item(126, 5).Value = 974
$ws.Cells.Item(126, 7).Value = 1
$ws.Cells.Item(126, 8).Value = 154

# Row 158
$ws.Cells.Item(158, 2).Value = 2305
$ws.Cells.Item(158, 3).Value = 11
$ws.Cells.Item(158, 4).Value = 1581
$ws.Cells.Item(158, 5).Value = 661

# Row 161
$ws.Cells.Item(161, 2).Value = 1972
$ws.Cells.Item(161, 3).Value = 23
$ws.Cells.Item(161, 4).Value = 1465
$ws.Cells.Item(161, 5).Value = 457
$ws.Cells.Item(161, 7).Value = 1
$ws.Cells.Item(161, 8).Value = 50

# Row 191
$ws.Cells.Item(191, 2).Value = 210
$ws.Cells.Item(191, 3).Value = 2
$ws.Cells.Item(191, 4).Value = 187
$ws.Cells.Item(191, 5).Value = 16
